$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_G_acc_LT"

$values = @(
    73.693693693693689,
    73.333333333333329,
    75.675675675675677,
    69.729729729729726,
    72.612612612612608,
    72.612612612612608,
    80.540540540540533,
    81.081081081081081,
    79.459459459459453,
    80.36036036036036,
    67.567567567567565,
    72.612612612612608,
    76.396396396396398,
    76.396396396396398,
    76.936936936936945,
    64.684684684684697,
    67.747747747747738,
    71.171171171171167,
    77.837837837837839,
    79.099099099099107,
    78.378378378378372,
    66.666666666666657,
    64.504504504504496,
    64.324324324324323,
    75.49549549549549,
    72.252252252252248,
    74.054054054054049,
    72.972972972972968,
    69.549549549549553,
    72.972972972972968,
    75.49549549549549,
    75.855855855855864,
    76.396396396396398,
    69.009009009009006,
    69.549549549549553,
    63.063063063063062,
    70.990990990990994,
    66.306306306306311,
    66.846846846846844,
    74.414414414414409,
    76.396396396396398,
    76.756756756756758,
    75.675675675675677,
    76.036036036036037,
    75.675675675675677,
    68.648648648648646,
    68.468468468468473,
    73.693693693693689
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
